$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 78
$ws1.Range("F3").Value = 204
$ws1.Range("F4").Value = 90
$ws1.Range("F5").Value = 1681
$ws1.Range("F6").Value = 3272
$ws1.Range("F7").Value = 881
$ws1.Range("F8").Value = 2085
$ws1.Range("F9").Value = 2000
$ws1.Range("F10").Value = 1033
$ws1.Range("F11").Value = 360
$ws1.Range("F13").Value = 1621
$ws1.Range("F14").Value = 350
$ws1.Range("F16").Value = 16
$ws1.Range("F18").Value = 89
$ws1.Range("F19").Value = 1463
$ws1.Range("F20").Value = 541
$ws1.Range("F21").Value = 648
$ws1.Range("F22").Value = 334
$ws1.Range("F23").Value = 11687
$ws1.Range("F24").Value = 11817
$ws1.Range("F25").Value = 866
$ws1.Range("F26").Value = 668
$ws1.Range("F27").Value = 1858
$ws1.Range("F28").Value = 161
$ws1.Range("F29").Value = 465

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 64

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 78
$ws4.Range("F3").Value = 64
$ws4.Range("F4").Value = 204
$ws4.Range("F6").Value = 90
$ws4.Range("F7").Value = 1681
$ws4.Range("F8").Value = 3272
$ws4.Range("F9").Value = 881
$ws4.Range("F10").Value = 2085
$ws4.Range("F11").Value = 2000
$ws4.Range("F12").Value = 1033
$ws4.Range("F13").Value = 360
$ws4.Range("F15").Value = 1621
$ws4.Range("F16").Value = 350
$ws4.Range("F18").Value = 16
$ws4.Range("F22").Value = 89
$ws4.Range("F23").Value = 1463
$ws4.Range("F24").Value = 541
$ws4.Range("F25").Value = 648
$ws4.Range("F26").Value = 334
$ws4.Range("F27").Value = 11689
$ws4.Range("F28").Value = 11817
$ws4.Range("F29").Value = 866
$ws4.Range("F30").Value = 668
$ws4.Range("F31").Value = 1858
$ws4.Range("F34").Value = 161
$ws4.Range("F35").Value = 465
